$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26). All rows below shift up by one.
$ws.Rows(26).Delete()

# Remove the "SC 92" row. After the first delete, it now sits at row 27.
$ws.Rows(27).Delete()

# --- Per-cell value corrections, expressed against the final (post-delete) row numbers ---

# Row 3 (RM 8): D3 gains a value
$ws.Range("D3").Value = -14.2

# Row 4 (RM 9): E4 becomes blank
$ws.Range("E4").ClearContents()

# Row 5 (RM 14): D5 becomes blank
$ws.Range("D5").ClearContents()

# Row 6 (RM 21): F6 gains a value
$ws.Range("F6").Value = 16.43

# Row 9 (RM 42): E9 gains a value
$ws.Range("E9").Value = -6.8

# Row 10 (RM 52 a): E10 gains a value
$ws.Range("E10").Value = -6.1

# Row 12 (RM 81): F12 becomes blank
$ws.Range("F12").ClearContents()

# Row 14 (RM 90): F14 gains a value
$ws.Range("F14").Value = 17.76

# Row 17 (RM 116): E17 becomes blank, F17 gains a value
$ws.Range("E17").ClearContents()
$ws.Range("F17").Value = 17.78

# Row 18 (RM 120): E18 becomes blank
$ws.Range("E18").ClearContents()

# Row 19 (RM 125): F19 gains a value
$ws.Range("F19").Value = 17.81

# Row 20 (RM 134): F20 becomes blank
$ws.Range("F20").ClearContents()

# Row 21 (RM 135): D21 gains a value
$ws.Range("D21").Value = -14.3

# Row 23 (RM 140): D23 becomes blank, F23 becomes blank
$ws.Range("D23").ClearContents()
$ws.Range("F23").ClearContents()

# Row 25 (RM 145): F25 becomes blank
$ws.Range("F25").ClearContents()

# Row 27 (SC 101): F27 gains a value (was blank)
$ws.Range("F27").Value = 17

# Row 28 (SC 105): F28 gains a value (was blank)
$ws.Range("F28").Value = 17.44

# Row 32 (SC 193): D32 gains a value (was blank)
$ws.Range("D32").Value = -14.7

echo "done"
